$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 22 with the next day's data (one day after row 21 / A21)
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$prevDate = $ws.Range("A21").Value2
$ws.Range("A22").Value2 = $prevDate + 1

$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 3

$ws.Range("K22").Select()
